$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row changes
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Row 2
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 0.04476192591403197
$ws.Range("F2").Value = 0.04084540229737399
$ws.Range("G2").Value = 0.03250659174605396

# Row 3
$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 0.04096842921706784
$ws.Range("F3").Value = 0.03738382428129747
$ws.Range("G3").Value = 0.02988305984496106

# Row 4
$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 0.04080585631739253
$ws.Range("F4").Value = 0.03723547598407179
$ws.Range("G4").Value = 0.02988518693521401

# Row 5
$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 0.03603028124857128
$ws.Range("F5").Value = 0.03287774827454609
$ws.Range("G5").Value = 0.02655562031246779

# Row 6
$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 0.03844007312878753
$ws.Range("F6").Value = 0.03507669116608782
$ws.Range("G6").Value = 0.02731884324769415

# Row 7
$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 0.03782418245151516
$ws.Range("F7").Value = 0.0345146889293497
$ws.Range("G7").Value = 0.02650634514331364

# Row 8
$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 0.02851984745386078
$ws.Range("F8").Value = 0.02602445312451372
$ws.Range("G8").Value = 0.01930674095221131
